$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "By profile" (G) column values for rows 3-9 with the
# newly-uploaded figures. Formulas in column H (which compute the
# deviation F/G-100%) recalculate automatically, as do the SUM totals
# in row 10.
$ws.Range("G3").Formula = "=73"
$ws.Range("G4").Formula = "=401"
$ws.Range("G5").Formula = "=312"
$ws.Range("G6").Formula = "=272"
$ws.Range("G7").Formula = "=272"
$ws.Range("G8").Formula = "=166"
$ws.Range("G9").Formula = "=183"

# Move the active cell/selection to match where the author left off.
$ws.Range("I10").Select()
